$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while preserving it as TEXT (not auto-converted
# to a number by Excel), and then strip the temporary "Text" number-format
# so no extra formatting is left behind on the cell.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# ---- Row 2 (Id 111545328 -> 111545414) ----
$ws.Range("A2").Value = 111545414
Set-TextValue "I2" "9"
Set-TextValue "J2" "registreringar"
$ws.Range("P2").Value = "Orsa Viborg, glänta i mitten av skogsparti, Dlr"
$ws.Range("Q2").Value = 480487.2503558649
$ws.Range("R2").Value = 6772784.264016891

# ---- Row 4 (Id 111545323 -> 111543957) ----
$ws.Range("A4").Value = 111543957
$ws.Range("B4").Value = 57494
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 205992
$ws.Range("F4").Value = "Vattenfladdermus"
$ws.Range("G4").Value = "Myotis daubentonii"
$ws.Range("H4").Value = "(Kuhl, 1817)"
Set-TextValue "I4" "1"
Set-TextValue "J4" "registreringar"
$ws.Range("P4").Value = "Orsa Viborg, glänta i skogsparti, Dlr"
$ws.Range("Q4").Value = 480406.6045043401
$ws.Range("R4").Value = 6772745.04339793

# ---- Row 5 (Id 111545414 -> 111545323) ----
$ws.Range("A5").Value = 111545323
$ws.Range("B5").Value = 57487
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 205998
$ws.Range("F5").Value = "Nordfladdermus"
$ws.Range("G5").Value = "Eptesicus nilssonii"
$ws.Range("H5").Value = "(A.Keyserling & Blasius, 1839)"
Set-TextValue "I5" "2"
$ws.Range("J5").ClearContents()
$ws.Range("P5").Value = "Orsa Viborg, intill en grupp med hålträd, Dlr"
$ws.Range("Q5").Value = 480427.8053356989
$ws.Range("R5").Value = 6772811.198980245

# ---- Row 7 (Id 111543957 -> 111545328) ----
$ws.Range("A7").Value = 111545328
$ws.Range("J7").ClearContents()
$ws.Range("P7").Value = "Orsa Viborg, intill en grupp med hålträd, Dlr"
$ws.Range("Q7").Value = 480427.8053356989
$ws.Range("R7").Value = 6772811.198980245
